# Working on weather data table: insert two new weather stations
# (Jena, GM and Flyvestation, DA) into the Latitude table, keeping the
# existing rows sorted by latitude, and update the view state (zoom +
# selection) to match the author's final editing position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above the current row 6 (Hannover, GM) for the
#    "Jena, GM" station -- it sits between Erfurt and Hannover by latitude.
$ws.Rows(6).Insert()
$ws.Range("A6").Value = "Jena, GM"
$ws.Range("B6").Value = 155
$ws.Range("C6").Value = 50.9267
$ws.Range("D6").Value = 11.5842
$ws.Range("E6").Value = 4

# 2) Insert a new row above the current row 11 (Oslo, NO, after the first
#    insert shifted it down to row 11) for the "Flyvestation, DA" station.
$ws.Rows(11).Insert()
$ws.Range("A11").Value = "Flyvestation, DA"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 57.093
$ws.Range("D11").Value = 9.849
$ws.Range("E11").Value = 0

# 3) Update the view: zoom to 125% and move the selection to the new last
#    row (A12), matching the active cell left after the edits.
$excel.ActiveWindow.Zoom = 125
$ws.Range("A12").Select()
